$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: standalone D-column (Malekiat_Dolat_Dar_Sherkat) corrections ---
$ws.Cells.Item(56, 4).Value2 = 0
$ws.Cells.Item(112, 4).Value2 = 0
$ws.Cells.Item(165, 4).Value2 = 0.3516
$ws.Cells.Item(335, 4).Value2 = 0
$ws.Cells.Item(359, 4).Value2 = 0.1104
$ws.Cells.Item(442, 4).Value2 = 0
$ws.Cells.Item(479, 4).Value2 = 0
$ws.Cells.Item(540, 4).Value2 = 0.701
$ws.Cells.Item(574, 4).Value2 = 0.3224
$ws.Cells.Item(629, 4).Value2 = 0
$ws.Cells.Item(633, 4).Value2 = 0.7128
$ws.Cells.Item(639, 4).Value2 = 0
$ws.Cells.Item(675, 4).Value2 = 0.8582

# --- Part 2: rows 685-953 -- Nationalcode/Name columns shift up by one row
# (a row was inserted/removed upstream in the source list; B/C realign to the
#  next row while D, mostly 0, stays put except for the one override below)
for ($r = 685; $r -le 953; $r++) {
    $srcB = $ws.Cells.Item($r + 1, 2).Value2
    $srcC = $ws.Cells.Item($r + 1, 3).Value2
    $ws.Cells.Item($r, 2).Value2 = $srcB
    $ws.Cells.Item($r, 3).Value2 = $srcC
}

# --- Part 3: D705 correction (part of the same shifted block) ---
$ws.Cells.Item(705, 4).Value2 = 0

# --- Part 4: row 954 -- new trailing entry (Bank Melli) added to the list ---
$ws.Cells.Item(954, 2).Value2 = 10861677542
$ws.Cells.Item(954, 3).Value2 = "بانک ملی"
